# Generate Report for Handoff
#
# Rows 7, 9, 10, 11, 12, 14 on the "zh-cn" and "de-de" sheets just finished
# a fresh handoff cycle: their "Priority" column is stamped "ht" (handoff
# type) and the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
# timestamps advance a few seconds to the moment the new handoff report was
# generated.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 11, 12, 14)

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

foreach ($r in $rows) {
    # Overview sheet: "Latest HO Xliff Generate Date" (column G)
    $overview.Range("G$r").Value = "2016-08-12 04:21:59"

    # zh-cn sheet: "Priority" (column E) and "Latest Handoff Datetime" (column H)
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-12 04:21:53"

    # de-de sheet: "Priority" (column E) and "Latest Handoff Datetime" (column H)
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-12 04:21:59"
}
